$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(6)

$tr = $sh.TextFrame.TextRange
$tr.Text = "Latency 		: 99thPercentile Latency"
$tr.InsertAfter("`rThroughput		: Operation/sec (YCSB의 target 옵션을 통해 제한하면서 실험 진행)") | Out-Null
$tr.InsertAfter("`rRecord 		: 100,000") | Out-Null
$tr.InsertAfter("`rWorkload		: Read(50%), Write(50%)") | Out-Null

# Resize / reposition the textbox (values chosen so float32 COM round-trip
# truncates to the exact target EMU: x=963600 y=4939475 cx=7216800 cy=1046700).
# Set after the text edits because this shape has spAutoFit, which otherwise
# recomputes Height from the new text content.
$sh.Left   = 75.8740157480315
$sh.Top    = 388.93504337007874
$sh.Width  = 568.251984503937
$sh.Height = 82.41732483464567
